$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "alt" column (J): header styled like the other header cells (copy I1's
# format onto J1), then the per-row team values.
$ws.Range("I1").Copy($ws.Range("J1")) | Out-Null
$ws.Range("J1").Value = "alt"

$teams = @("team1","team2","team3","team4","team5","team6","team7","team8","team9","team10")
for ($i = 0; $i -lt $teams.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 10).Value = $teams[$i]
}

# Narrow column B to a manual (non bestfit) width. The target stored width is
# 8.1796875 characters; this engine quantizes ColumnWidth writes to steps of
# 1/6 character, so 7.3 is the closest input that lands on the nearest
# achievable stored width (8.16666...) to the target.
$ws.Columns("B").ColumnWidth = 7.3

# Reflect the new selection left by the editor
$ws.Range("I15").Select() | Out-Null

Write-Output "done"
